$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column A entirely, shifting the remaining columns (B:F) left into A:E
$ws.Columns.Item(1).Delete()

# Fix the header text: MODEL_CONDITION -> MODELCONDITION
$ws.Cells.Replace("MODEL_CONDITION", "MODELCONDITION")
